$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is unambiguously non-numeric (contains extra
# separators, subscript digits, percent signs, or padding) - can be
# written directly and Excel will keep them as text, matching the
# original inlineStr cells.
$ws.Range('D2').Value = '57.930.01'
$ws.Range('E2').Value = '  -1.81%  '
$ws.Range('D3').Value = '3.093.93'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('E6').Value = '  -2.32%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.093.86'
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('E10').Value = '  -2.99%  '
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('E12').Value = '  +2.52%  '
$ws.Range('D13').Value = '3.625.92'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('E14').Value = '  +2.39%  '
$ws.Range('E15').Value = '  -5.65%  '
$ws.Range('E16').Value = '  -1.42%  '
$ws.Range('D17').Value = '57.969.02'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').Value = '3.089.64'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('E20').Value = '  -2.84%  '
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('E25').Value = '  +2.56%  '
$ws.Range('E26').Value = '  -1.38%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '0.0₃0916'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('E30').Value = '  -5.03%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('E32').Value = '  +1.63%  '
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('E35').Value = '  +2.03%  '
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('E38').Value = '  -3.23%  '
$ws.Range('E39').Value = '  -5.34%  '
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('E41').Value = '  +1.31%  '
$ws.Range('E42').Value = '  +5.59%  '
$ws.Range('E43').Value = '  +2.61%  '
$ws.Range('D44').Value = '3.133.85'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('E47').Value = '  +2.10%  '
$ws.Range('D48').Value = '2.272.44'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('E49').Value = '  +2.72%  '
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('E51').Value = '  -2.76%  '

# Cells whose new text looks like a plain number (e.g. "7.16") would
# otherwise be auto-converted to a numeric value by Excel. Temporarily
# mark them as Text so the literal string is preserved, then restore
# the default "Normal" style so no stray formatting is left behind.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '525.96'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '140.74'
$ws.Range('D6').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.16'
$ws.Range('D10').Style = "Normal"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.392'
$ws.Range('D12').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '25.61'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000165'
$ws.Range('D16').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.09'
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.68'
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.95'
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '342.24'
$ws.Range('D22').Style = "Normal"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.511'
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '67.45'
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.170'
$ws.Range('D26').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.39'
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.26'
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.87'
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.95'
$ws.Range('D33').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '158.42'
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.62'
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.14'
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '26.09'
$ws.Range('D38').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0668'
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.01'
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.54'
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.684'
$ws.Range('D43').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '36.93'
$ws.Range('D45').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0263'
$ws.Range('D47').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.988'
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.09'
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '20.42'
$ws.Range('D51').Style = "Normal"
